# Publish documentation 0.1.1 / ror 0.1.1
# - Bump the "Version" value on the Metadata sheet from 0.1.0 to 0.1.1
# - Bump the "Date" value on the Metadata sheet to the new publication date
# - Add a new "Context" row for the new "element:Address" extension context

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Version (row 3) and Date (row 8) values
$ws.Range("B3").Value = "0.1.1"
$ws.Range("B8").Value = "2023-06-02T12:02:38+02:00"

# Append a new Context row (row 21) for "element:Address", copying the
# formatting of the existing Context row (row 20) and then setting the text.
$ws.Range("A20:B20").Copy()
$ws.Range("A21:B21").PasteSpecial(-4122)
$ws.Range("A21").Value = "Context"
$ws.Range("B21").Value = "element:Address"
